# feat: filtrar gmb y reportes de gd
#
# Updates row 2 (existing product) with new values and appends two new
# product rows (3 and 4) below it. Every cell in this sheet is plain
# text (inlineStr) even when it looks numeric, so any cell whose new
# value could be mis-read as a number is forced to Text format first -
# this keeps leading zeros (e.g. "000000600453") and exact decimal
# literals (e.g. "3.2030357142860004") intact instead of letting Excel
# silently reinterpret/round them as numbers. Cells that already hold
# (or keep holding) plain non-numeric text are left alone so we don't
# touch formatting that isn't actually changing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: only the columns that actually change -------------------
# (NumberFormat is applied per contiguous block; this engine's COM
# layer does not reliably honour multi-area "A1,B2:C2" range strings.)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("C2:D2").NumberFormat = "@"
$ws.Range("F2:J2").NumberFormat = "@"

$ws.Range("A2").Value  = "000000600453"
$ws.Range("B2").Value  = "COLADOR 12 EN 1"
$ws.Range("C2").Value  = "11"
$ws.Range("D2").Value  = "12"
$ws.Range("F2").Value  = "7.197500000000001"
$ws.Range("G2").Value  = "4.445178571429"
$ws.Range("H2").Value  = "3.552410714286"
$ws.Range("I2").Value  = "3.347053571429"
$ws.Range("J2").Value  = "3.2030357142860004"
$ws.Range("M2").Value  = "COCINA"
$ws.Range("AC2").Value = "HOCO"

# --- New rows 3 and 4 -------------------------------------------------
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE")

$row3 = @{
    A="009283044565"; B="VENDAS DE BOX"; C="207"; D="207"; E="";
    F="5"; G="3.9"; H="3.75"; I="3.25"; J="2.95"; K="0"; L="";
    M="ACCESORIOS PARA EJERCICIOS"; N="NORMAL"; O="1"; P="0"; Q="0";
    R="IVA 15%"; S="15"; T="NO"; U=""; V=""; W=""; X=""; Y=""; Z=""; AA="";
    AB="NO"; AC="DEEJ"; AD="DE"; AE="DEPORTES"
}

$row4 = @{
    A="17874117802"; B="DISPENSADOR DE JUGO X3"; C="6"; D="6"; E="";
    F="34.25276785714001"; G="31.15160714286"; H="27.40410714286";
    I="24.895552321430003"; J="23.803571741069998"; K="0"; L="";
    M="ACCESORIOS PARA EL HOGAR"; N="NORMAL"; O="1"; P="0"; Q="0";
    R="IVA 15%"; S="15"; T="NO"; U=""; V=""; W=""; X=""; Y=""; Z=""; AA="";
    AB="NO"; AC="HOAC"; AD="HO"; AE="HOGAR"
}

$rows = @{ 3 = $row3; 4 = $row4 }

foreach ($r in 3,4) {
    $data = $rows[$r]
    $ws.Range("A$r`:AE$r").NumberFormat = "@"
    foreach ($col in $columns) {
        $ws.Range("$col$r").Value = $data[$col]
    }
}
